$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the base value which drives the dependent formulas
$ws.Range("A1").Value = 1

# Update the active selection to A2
$ws.Range("A2").Select()
